# Simplify unnecessarily long PV names
#
# Renames the "PROP" (column D) shared-string values used to build the
# PV NAME column (F = C & D & E) for the virtual/physical actual
# position/velocity rows (rows 29-44) on every module sheet.

$wb = $excel.ActiveWorkbook

# Old name -> new (simplified) name.
$renameMap = @{
    "VirtCIDActualPos"  = "CIDVirtPos"
    "VirtCIEActualPos"  = "CIEVirtPos"
    "VirtCSDActualPos"  = "CSDVirtPos"
    "VirtCSEActualPos"  = "CSEVirtPos"
    "VirtCIDActualVelo" = "CIDVirtVelo"
    "VirtCIEActualVelo" = "CIEVirtVelo"
    "VirtCSDActualVelo" = "CSDVirtVelo"
    "VirtCSEActualVelo" = "CSEVirtVelo"
    "PhyCIDActualPos"   = "CIDPhyPos"
    "PhyCIEActualPos"   = "CIEPhyPos"
    "PhyCSDActualPos"   = "CSDPhyPos"
    "PhyCSEActualPos"   = "CSEPhyPos"
    "PhyCIDActualVelo"  = "CIDPhyVelo"
    "PhyCIEActualVelo"  = "CIEPhyVelo"
    "PhyCSDActualVelo"  = "CSDPhyVelo"
    "PhyCSEActualVelo"  = "CSEPhyVelo"
}

foreach ($sheetName in @("Mod01", "Mod02", "Mod03")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($row = 29; $row -le 44; $row++) {
        $cell = $ws.Cells.Item($row, 4)   # column D = "PROP"
        $old = $cell.Value2
        if ($renameMap.ContainsKey($old)) {
            $cell.Value = $renameMap[$old]
        }
    }
}

# Scroll the sheet-tab strip back so "Mod01" is the first displayed tab
# (previously it was scrolled to show "Mod02" first).
$excel.ActiveWindow.DisplayedFirstSheet = "Mod01"
